# First pass at power calculations using simulation
# Adds a new "Time" column (H) with "Before"/"After" values for each phrase.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("H1").Value = "Time"
$ws.Range("H1").Font.Bold = $true

# Per-row values for the new "Time" column
$ws.Range("H2").Value = "Before"
$ws.Range("H3").Value = "After"
$ws.Range("H4").Value = "After"
$ws.Range("H5").Value = "Before"
$ws.Range("H6").Value = "After"
$ws.Range("H7").Value = "Before"
$ws.Range("H8").Value = "After"
$ws.Range("H9").Value = "After"
$ws.Range("H10").Value = "Before"
$ws.Range("H11").Value = "After"
$ws.Range("H12").Value = "After"
$ws.Range("H13").Value = "Before"
$ws.Range("H14").Value = "After"
$ws.Range("H15").Value = "Before"
$ws.Range("H16").Value = "Before"
$ws.Range("H17").Value = "After"

# Reflect the selection change recorded in the workbook view
$null = $ws.Range("J18").Select()
